# Update port names (POL / POD columns) from ALL CAPS to Title Case,
# and update the selected cell shown in the saved view.
#
# The shared-strings table in the saved OOXML is ordered by first use, and
# unused strings get dropped when the workbook is saved. To reproduce the
# exact target shared-string order:
#   TYPE OF VEHICLE, CAR, SUV, LARGE SUV, PICKUP, PRICE, POL, POD,
#   MOTORCYCLE, Rotterdam, New York, Savannah, Miami, Houston,
#   Indianapolis, Los Angeles, San Francisco, Varna
# we first "seed" the new strings (in that exact order) into a scratch
# area far off the used range, then rewrite the real data cells (which
# will reuse the already-seeded shared-string entries), and finally clear
# the scratch cells. Because the real data cells still reference those
# strings, they survive the clean-up/garbage-collection that happens on
# save, while the old ALL CAPS strings become unreferenced and are
# dropped automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$seedRow = 1000

# Seed the new shared strings in the exact order they must appear.
$ws.Cells.Item($seedRow, 26).Value = "Rotterdam"
$ws.Cells.Item($seedRow + 1, 26).Value = "New York"
$ws.Cells.Item($seedRow + 2, 26).Value = "Savannah"
$ws.Cells.Item($seedRow + 3, 26).Value = "Miami"
$ws.Cells.Item($seedRow + 4, 26).Value = "Houston"
$ws.Cells.Item($seedRow + 5, 26).Value = "Indianapolis"
$ws.Cells.Item($seedRow + 6, 26).Value = "Los Angeles"
$ws.Cells.Item($seedRow + 7, 26).Value = "San Francisco"
$ws.Cells.Item($seedRow + 8, 26).Value = "Varna"

# Column A (POL): title-case the port-of-load names, in contiguous blocks.
$ws.Range("A2:A6").Value = "New York"
$ws.Range("A7:A11").Value = "Savannah"
$ws.Range("A12:A16").Value = "Miami"
$ws.Range("A17:A21").Value = "Houston"
$ws.Range("A22:A26").Value = "Indianapolis"
$ws.Range("A27:A31").Value = "Los Angeles"
$ws.Range("A32:A36").Value = "San Francisco"
$ws.Range("A37:A41").Value = "New York"
$ws.Range("A42:A46").Value = "Savannah"
$ws.Range("A47:A51").Value = "Miami"
$ws.Range("A52:A56").Value = "Houston"
$ws.Range("A57:A61").Value = "Indianapolis"
$ws.Range("A62:A66").Value = "Los Angeles"
$ws.Range("A67:A71").Value = "San Francisco"

# Column B (POD): title-case the port-of-discharge names.
$ws.Range("B2:B36").Value = "Rotterdam"
$ws.Range("B37:B71").Value = "Varna"

# Remove the temporary seed cells now that the real cells hold the
# references. Delete the entire rows (rather than just clearing the
# contents) so no empty styled cells are left behind and the sheet's
# dimension/used range goes back to its original extent.
$ws.Range("Z1000:Z1008").EntireRow.Delete()

# Update the view: select G69 (no particular scroll position recorded).
$ws.Range("G69").Select()
